# Updated symbol list (crypto price/volume refresh) — mirrors the
# GitHub Actions scraper commit. Each cell in columns D (Price) and
# E (Volume(1h)) holds a plain text value (e.g. "304.43", "-3.81%").
# We prefix the literal with a leading apostrophe so Excel stores it
# as text (quote-prefixed), exactly like the source data, instead of
# auto-converting the number-looking / percent-looking strings into
# numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("E2").Value = "'-0.31%"

$ws.Range("D3").Value = "'35.55"
$ws.Range("E3").Value = "'-4.05%"

$ws.Range("D4").Value = "'5.060"
$ws.Range("E4").Value = "'0.98%"

$ws.Range("D5").Value = "'0.07888"
$ws.Range("E5").Value = "'-0.04%"

$ws.Range("D6").Value = "'2.104"
$ws.Range("E6").Value = "'-4.85%"

$ws.Range("D7").Value = "'4.120"
$ws.Range("E7").Value = "'2.69%"

$ws.Range("D8").Value = "'7.907"
$ws.Range("E8").Value = "'-1.34%"

$ws.Range("D9").Value = "'0.9222"
$ws.Range("E9").Value = "'0.23%"

$ws.Range("D10").Value = "'0.09743"
$ws.Range("E10").Value = "'0.96%"

$ws.Range("D11").Value = "'0.1840"
$ws.Range("E11").Value = "'-2.64%"

$ws.Range("D12").Value = "'0.08654"
$ws.Range("E12").Value = "'0.93%"

$ws.Range("D13").Value = "'0.03570"
$ws.Range("E13").Value = "'-3.03%"

$ws.Range("D14").Value = "'0.09915"
$ws.Range("E14").Value = "'-0.71%"

$ws.Range("D15").Value = "'0.001431"
$ws.Range("E15").Value = "'-3.23%"

$ws.Range("D16").Value = "'0.005662"
$ws.Range("E16").Value = "'0.36%"

$ws.Range("D17").Value = "'3.471"
$ws.Range("E17").Value = "'0.22%"

$ws.Range("D18").Value = "'2.641"
$ws.Range("E18").Value = "'17.12%"

$ws.Range("D19").Value = "'0.3370"
$ws.Range("E19").Value = "'-1.29%"

$ws.Range("D20").Value = "'0.1340"
$ws.Range("E20").Value = "'1.83%"

$ws.Range("D21").Value = "'5.177"
$ws.Range("E21").Value = "'8.83%"

$ws.Range("E22").Value = "'0.61%"

$ws.Range("D23").Value = "'0.04496"
$ws.Range("E23").Value = "'-1.06%"

$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'0.00%"

$ws.Range("D25").Value = "'0.004856"
$ws.Range("E25").Value = "'8.57%"

$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-6.92%"

$ws.Range("E27").Value = "'0.15%"

$ws.Range("D39").Value = "'0.01835"
$ws.Range("E39").Value = "'-0.66%"

$ws.Range("D40").Value = "'0.04718"
$ws.Range("E40").Value = "'-0.99%"

$ws.Range("D41").Value = "'0.007905"
$ws.Range("E41").Value = "'-2.66%"

$ws.Range("D42").Value = "'0.1389"
$ws.Range("E42").Value = "'-0.81%"

$ws.Range("E43").Value = "'2.56%"

$ws.Range("D44").Value = "'0.002193"
$ws.Range("E44").Value = "'-1.60%"

$ws.Range("D45").Value = "'0.01117"
$ws.Range("E45").Value = "'5.84%"

$ws.Range("D46").Value = "'0.00006282"
$ws.Range("E46").Value = "'-0.14%"

$ws.Range("E48").Value = "'0.28%"

$ws.Range("D49").Value = "'50.68"
$ws.Range("E49").Value = "'69.06%"

$ws.Range("E50").Value = "'10.70%"
